$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the "mobile network coverage by technology" table with a new
# "2022" column (N), mirroring the formatting of the existing "2021"
# column (M) for the header row and the 2G/3G/4G data rows.
$ws.Range("M3:M7").Copy($ws.Range("N3:N7"))

$ws.Range("N4").Value = 2022
$ws.Range("N5").Value = 98.8
$ws.Range("N6").Value = 98
$ws.Range("N7").Value = 96.9

# The workbook was left with cell O4 selected.
$ws.Range("O4").Select() | Out-Null
